$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'English'
$ws.Range("B1").Value = 'Japanese'
$ws.Range("A2").Value = 'internet'
$ws.Range("B2").Value = 'インターネット'
$ws.Range("A3").Value = 'painting; picture; drawing'
$ws.Range("B3").Value = '絵|え'
$ws.Range("A4").Value = 'movie theater'
$ws.Range("B4").Value = '映画館|えいがかん'
$ws.Range("A5").Value = 'foreigner'
$ws.Range("B5").Value = '外国人|がいこくじん'
$ws.Range("A6").Value = 'furniture'
$ws.Range("B6").Value = '家具|かぐ'
$ws.Range("A7").Value = 'night club'
$ws.Range("B7").Value = 'クラブ'
$ws.Range("A8").Value = 'experience'
$ws.Range("B8").Value = '経験|けいけん'
$ws.Range("A9").Value = 'wedding'
$ws.Range("B9").Value = '結婚式|けっこんしき'
$ws.Range("A10").Value = 'earthquake'
$ws.Range("B10").Value = '地震|じしん'
$ws.Range("A11").Value = 'deadline'
$ws.Range("B11").Value = '締め切り|しめきり'
$ws.Range("A12").Value = 'jacket'
$ws.Range("B12").Value = 'ジャケット'
$ws.Range("A13").Value = 'custom'
$ws.Range("B13").Value = '習慣|しゅうかん'
$ws.Range("A14").Value = 'tax'
$ws.Range("B14").Value = '税金|ぜいきん'
$ws.Range("A15").Value = 'graduation ceremony'
$ws.Range("B15").Value = '卒業式|そつぎょうしき'
$ws.Range("A16").Value = 'soba; Japanese buckwheat noodles'
$ws.Range("B16").Value = 'そば'
$ws.Range("A17").Value = 'map'
$ws.Range("B17").Value = '地図|ちず'
$ws.Range("A18").Value = 'battery'
$ws.Range("B18").Value = '電池|でんち'
$ws.Range("A19").Value = 'garden'
$ws.Range("B19").Value = '庭|にわ'
$ws.Range("A20").Value = 'presentation'
$ws.Range("B20").Value = '発表|はっぴょう'
$ws.Range("A21").Value = 'broadcast program'
$ws.Range("B21").Value = '番組|ばんぐみ'
$ws.Range("A22").Value = 'swimming pool'
$ws.Range("B22").Value = 'プール'
$ws.Range("A23").Value = 'pet'
$ws.Range("B23").Value = 'ペット'
$ws.Range("A24").Value = 'insurance'
$ws.Range("B24").Value = '保険|ほけん'
$ws.Range("A25").Value = 'schedule; plan'
$ws.Range("B25").Value = '予定|よてい'
$ws.Range("A26").Value = 'Japanese inn'
$ws.Range("B26").Value = '旅館|りょかん'
$ws.Range("A27").Value = 'discount coupon'
$ws.Range("B27").Value = '割引券|わりびきけん'
$ws.Range("A28").Value = 'spacious; wide'
$ws.Range("B28").Value = '広い|ひろい'
$ws.Range("A29").Value = 'to sell'
$ws.Range("B29").Value = '売る|うる'
$ws.Range("A30").Value = 'to withdraw (money)'
$ws.Range("B30").Value = '下ろす|おろす'
$ws.Range("A31").Value = 'to draw; to paint'
$ws.Range("B31").Value = '描く|かく'
$ws.Range("A32").Value = 'to look for'
$ws.Range("B32").Value = '探す|さがす'
$ws.Range("A33").Value = 'to invite'
$ws.Range("B33").Value = '誘う|さそう'
$ws.Range("A34").Value = 'to chat'
$ws.Range("B34").Value = 'しゃべる'
$ws.Range("A35").Value = 'to date (someone); to keep company'
$ws.Range("B35").Value = '付き合う|つきあう'
$ws.Range("A36").Value = 'to arrive'
$ws.Range("B36").Value = '着く|つく'
$ws.Range("A37").Value = 'to buy insurance'
$ws.Range("B37").Value = '保険に入る|ほけんにはいる'
$ws.Range("A38").Value = 'to be cautious/careful'
$ws.Range("B38").Value = '気をつける|きをつける'
$ws.Range("A39").Value = 'to look into (a matter)'
$ws.Range("B39").Value = '調べる|しらべる'
$ws.Range("A40").Value = 'to be visible'
$ws.Range("B40").Value = '見える|みえる'
$ws.Range("A41").Value = 'to do sightseeing'
$ws.Range("B41").Value = '観光する|かんこうする'
$ws.Range("A42").Value = 'to decide on (an item)'
$ws.Range("B42").Value = 'する'
$ws.Range("A43").Value = 'to graduate (from...)'
$ws.Range("B43").Value = '卒業する|そつぎょうする'
$ws.Range("A44").Value = 'to reserve'
$ws.Range("B44").Value = '予約する|よやくする'
$ws.Range("A45").Value = 'all day long'
$ws.Range("B45").Value = '一日中|いちにちじゅう'
$ws.Range("A46").Value = '...,but; ...,so'
$ws.Range("B46").Value = '～けど'
$ws.Range("A47").Value = 'recently'
$ws.Range("B47").Value = '最近|さいきん'
$ws.Range("A48").Value = 'number...'
$ws.Range("B48").Value = '～番|～ばん'
$ws.Range("A49").Value = '-th'
$ws.Range("B49").Value = '～目|～め'
$ws.Range("A50").Value = 'one more time'
$ws.Range("B50").Value = 'もう一度|もういちど'
$ws.Range("A51").Value = '...nights'
$ws.Range("B51").Value = '～泊|～はく'
$ws.Range("A52").Value = 'with...'
$ws.Range("B52").Value = '～付|～つき'
$ws.Range("A53").Value = 'with meals'
$ws.Range("B53").Value = '食事付|しょくじつき'
$ws.Range("A54").Value = 'one night with two meals'
$ws.Range("B54").Value = '一泊二食付|いっぱくにしょくつき'
$ws.Range("A55").Value = 'checking in'
$ws.Range("B55").Value = 'チェックイン（する）'
$ws.Range("A56").Value = 'checking out'
$ws.Range("B56").Value = 'チェックアウト（する）'
$ws.Range("A57").Value = 'single room'
$ws.Range("B57").Value = 'シングル'
$ws.Range("A58").Value = 'double room'
$ws.Range("B58").Value = 'ダブル'
$ws.Range("A59").Value = 'twin room'
$ws.Range("B59").Value = 'ツイン'
$ws.Range("A60").Value = '...person(s)'
$ws.Range("B60").Value = '～名|～めい'
$ws.Range("A61").Value = 'receptionist; front desk'
$ws.Range("B61").Value = 'フロント'
$ws.Range("A62").Value = 'non-smoking room'
$ws.Range("B62").Value = '禁煙ルーム|きんえんルーム'
$ws.Range("A63").Value = 'smoking room'
$ws.Range("B63").Value = '喫煙ルーム|きつえんルーム'
$ws.Range("A64").Value = 'Can I pay by credit card?'
$ws.Range("B64").Value = 'クレジットカードで払えますか。|クレジットカードではらえますか。'
$ws.Range("A65").Value = 'Could you keep my luggage until 2 o''clock?'
$ws.Range("B65").Value = '二時まで荷物を預かってくれませんか。|にじまでにもつをあずかってくれませんか。'
$ws.Range("A66").Value = 'first person'
$ws.Range("B66").Value = '一人目|ひとりめ'
$ws.Range("A67").Value = 'second person'
$ws.Range("B67").Value = '二人目|ふたりめ'
$ws.Range("A68").Value = 'third person'
$ws.Range("B68").Value = '三人目|さんにんめ'
$ws.Range("A69").Value = 'first sheet'
$ws.Range("B69").Value = '一枚目|いちまいめ'
$ws.Range("A70").Value = 'second sheet'
$ws.Range("B70").Value = '二枚目|にまいめ'
$ws.Range("A71").Value = 'third sheet'
$ws.Range("B71").Value = '三枚目|さんまいめ'
$ws.Range("A72").Value = 'first year'
$ws.Range("B72").Value = '一年目|いちねんめ'
$ws.Range("A73").Value = 'second year'
$ws.Range("B73").Value = '二年目|にねんめ'
$ws.Range("A74").Value = 'third year'
$ws.Range("B74").Value = '三年目|さんねんめ'
$ws.Range("A75").Value = 'first day'
$ws.Range("B75").Value = '一日目|いちにちめ'
$ws.Range("A76").Value = 'second day'
$ws.Range("B76").Value = '二日目|ふつかめ'
$ws.Range("A77").Value = 'third day'
$ws.Range("B77").Value = '三日目|みっかめ'
$ws.Range("A78").Value = 'to die'
$ws.Range("B78").Value = '死ぬ|しぬ'
$ws.Range("A79").Value = 'death'
$ws.Range("B79").Value = '死|し'
$ws.Range("A80").Value = 'desperate'
$ws.Range("B80").Value = '必死|ひっし'
$ws.Range("A81").Value = 'the dead'
$ws.Range("B81").Value = '死者|ししゃ'
$ws.Range("A82").Value = 'meaning'
$ws.Range("B82").Value = '意味|いみ'
$ws.Range("A83").Value = 'to watch out'
$ws.Range("B83").Value = '注意する|ちゅういする'
$ws.Range("A84").Value = 'opinion'
$ws.Range("B84").Value = '意見|いけん'
$ws.Range("A85").Value = 'to prepare'
$ws.Range("B85").Value = '用意する|よういする'
$ws.Range("A86").Value = 'hobby'
$ws.Range("B86").Value = '趣味|しゅみ'
$ws.Range("A87").Value = 'interest'
$ws.Range("B87").Value = '興味|きょうみ'
$ws.Range("A88").Value = 'soybean paste'
$ws.Range("B88").Value = '味噌|みそ'
$ws.Range("A89").Value = 'taste'
$ws.Range("B89").Value = '味|あじ'
$ws.Range("A90").Value = 'to watch out'
$ws.Range("B90").Value = '注意する|ちゅういする'
$ws.Range("A91").Value = 'to order'
$ws.Range("B91").Value = '注文する|ちゅうもんする'
$ws.Range("A92").Value = 'to pour'
$ws.Range("B92").Value = '注ぐ|そそぐ'
$ws.Range("A93").Value = 'summer'
$ws.Range("B93").Value = '夏|なつ'
$ws.Range("A94").Value = 'summer vacation'
$ws.Range("B94").Value = '夏休み|なつやすみ'
$ws.Range("A95").Value = 'early summer'
$ws.Range("B95").Value = '初夏|しょか'
$ws.Range("A96").Value = 'fish'
$ws.Range("B96").Value = '魚|さかな'
$ws.Range("A97").Value = 'fish market'
$ws.Range("B97").Value = '魚市場|うおいちば'
$ws.Range("A98").Value = 'gold fish'
$ws.Range("B98").Value = '金魚|きんぎょ'
$ws.Range("A99").Value = 'mermaid'
$ws.Range("B99").Value = '人魚|にんぎょ'
$ws.Range("A100").Value = 'temple'
$ws.Range("B100").Value = 'お寺|おてら'
$ws.Range("A101").Value = 'Toji Temple'
$ws.Range("B101").Value = '東寺|とうじ'
$ws.Range("A102").Value = 'sacred building'
$ws.Range("B102").Value = '寺院|じいん'
$ws.Range("A103").Value = 'zen temple'
$ws.Range("B103").Value = '禅寺|ぜんでら'
$ws.Range("A104").Value = 'wide; spacious'
$ws.Range("B104").Value = '広い|ひろい'
$ws.Range("A105").Value = 'square; open space'
$ws.Range("B105").Value = '広場|ひろば'
$ws.Range("A106").Value = 'Hiroshima'
$ws.Range("B106").Value = '広島|ひろしま'
$ws.Range("A107").Value = 'advertisement'
$ws.Range("B107").Value = '広告|こうこく'
$ws.Range("A108").Value = 'bicycle'
$ws.Range("B108").Value = '自転車|じてんしゃ'
$ws.Range("A109").Value = 'to drive'
$ws.Range("B109").Value = '運転する|うんてんする'
$ws.Range("A110").Value = 'rotating sushi'
$ws.Range("B110").Value = '回転ずし|かいてんずし'
$ws.Range("A111").Value = 'to tumble; to fall down'
$ws.Range("B111").Value = '転ぶ|ころぶ'
$ws.Range("A112").Value = 'to borrow'
$ws.Range("B112").Value = '借りる|かりる'
$ws.Range("A113").Value = 'rented land'
$ws.Range("B113").Value = '借地|しゃくち'
$ws.Range("A114").Value = 'debt'
$ws.Range("B114").Value = '借金|しゃっきん'
$ws.Range("A115").Value = 'rented house'
$ws.Range("B115").Value = '借家|しゃくや'
$ws.Range("A116").Value = 'to run'
$ws.Range("B116").Value = '走る|はしる'
$ws.Range("A117").Value = 'hasty writing'
$ws.Range("B117").Value = '走り書き|はしりがき'
$ws.Range("A118").Value = 'escape from a prison'
$ws.Range("B118").Value = '脱走|だっそう'
$ws.Range("A119").Value = 'building'
$ws.Range("B119").Value = '建物|たてもの'
$ws.Range("A120").Value = 'to build'
$ws.Range("B120").Value = '建てる|たてる'
$ws.Range("A121").Value = 'to be built'
$ws.Range("B121").Value = '建つ|たつ'
$ws.Range("A122").Value = 'founding a nation'
$ws.Range("B122").Value = '建国|けんこく'
$ws.Range("A123").Value = 'underground'
$ws.Range("B123").Value = '地下|ちか'
$ws.Range("A124").Value = 'subway'
$ws.Range("B124").Value = '地下鉄|ちかてつ'
$ws.Range("A125").Value = 'map'
$ws.Range("B125").Value = '地図|ちず'
$ws.Range("A126").Value = 'earth; globe'
$ws.Range("B126").Value = '地球|ちきゅう'
$ws.Range("A127").Value = 'earthquake'
$ws.Range("B127").Value = '地震|じしん'
$ws.Range("A128").Value = 'square; open space'
$ws.Range("B128").Value = '広場|ひろば'
$ws.Range("A129").Value = 'place'
$ws.Range("B129").Value = '場所|ばしょ'
$ws.Range("A130").Value = 'case'
$ws.Range("B130").Value = '場合|ばあい'
$ws.Range("A131").Value = 'parking garage'
$ws.Range("B131").Value = '駐車場|ちゅうしゃじょう'
$ws.Range("A132").Value = 'foot; leg'
$ws.Range("B132").Value = '足|あし'
$ws.Range("A133").Value = 'to be sufficient'
$ws.Range("B133").Value = '足りる|たりる'
$ws.Range("A134").Value = 'one pair of shoes'
$ws.Range("B134").Value = '一足|いっそく'
$ws.Range("A135").Value = 'lack of water'
$ws.Range("B135").Value = '水不足|みずぶそく'
$ws.Range("A136").Value = 'to go through; to pass'
$ws.Range("B136").Value = '通る|とおる'
$ws.Range("A137").Value = 'to commute'
$ws.Range("B137").Value = '通う|かよう'
$ws.Range("A138").Value = 'going to school'
$ws.Range("B138").Value = '通学|つうがく'
$ws.Range("A139").Value = 'going to work'
$ws.Range("B139").Value = '通勤|つうきん'
